$wb = $excel.ActiveWorkbook

$newFile = "1d7d6837-711c-4007-a6ec-6372acb2f548.md"
$newStatus = "Handoff transform failed"
$epoch = "0001-01-01 00:00:00"
$ignored = "Ignored"

# ---------- Overview sheet ----------
$wsO = $wb.Worksheets.Item("Overview")
$wsO.Range("B2").Value = $newStatus
$wsO.Range("C2").Value = $newStatus

$wsO.Cells.Hyperlinks.Delete()
$oMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/23d5b9c3ebe74f1ef7e3b774a2e95a01f5f7c52d/e2e/" + $newFile
$oCfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/23d5b9c3ebe74f1ef7e3b774a2e95a01f5f7c52d/.localization-config"
$wsO.Hyperlinks.Add($wsO.Range("A2"), $oMdUrl, "", "", $newFile)
$wsO.Hyperlinks.Add($wsO.Range("A3"), $oCfgUrl, "", "", ".localization-config")

# ---------- zh-cn sheet ----------
$ws1 = $wb.Worksheets.Item("zh-cn")
$ws1.Range("B2").Value = $newStatus
$ws1.Range("C2").Clear()
$ws1.Range("D2").Value = $epoch
$ws1.Range("G2").Value = $epoch
$ws1.Range("H2").Value = $ignored

$ws1.Cells.Hyperlinks.Delete()
$z1MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/23d5b9c3ebe74f1ef7e3b774a2e95a01f5f7c52d/e2e/" + $newFile
$z1CfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/23d5b9c3ebe74f1ef7e3b774a2e95a01f5f7c52d/.localization-config"
$ws1.Hyperlinks.Add($ws1.Range("A2"), $z1MdUrl, "", "", $newFile)
$ws1.Hyperlinks.Add($ws1.Range("A3"), $z1CfgUrl, "", "", ".localization-config")

# ---------- de-de sheet ----------
$ws2 = $wb.Worksheets.Item("de-de")
$ws2.Range("B2").Value = $newStatus
$ws2.Range("C2").Clear()
$ws2.Range("D2").Value = $epoch
$ws2.Range("G2").Value = $epoch
$ws2.Range("H2").Value = $ignored

$ws2.Cells.Hyperlinks.Delete()
$d1MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/23d5b9c3ebe74f1ef7e3b774a2e95a01f5f7c52d/e2e/" + $newFile
$d1CfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/23d5b9c3ebe74f1ef7e3b774a2e95a01f5f7c52d/.localization-config"
$ws2.Hyperlinks.Add($ws2.Range("A2"), $d1MdUrl, "", "", $newFile)
$ws2.Hyperlinks.Add($ws2.Range("A3"), $d1CfgUrl, "", "", ".localization-config")
